# Add new "2025" worksheet for Policia data (commit: "added new data for Policia")

$wb = $excel.ActiveWorkbook

# --- Update the selection on the previously active sheet ("2024") ---
# Before switching away from it, the author had selected the whole data range (A1:C14).
$ws2024 = $wb.Worksheets.Item("2024")
$ws2024.Select() | Out-Null
$ws2024.Range("A1:C14").Select() | Out-Null

# --- Create the new "2025" sheet, placed after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025"

# Header row
$ws.Range("A1").Value = "Delitos"
$ws.Range("B1").Value = "Mujeres"
$ws.Range("C1").Value = "Hombres"

# Data rows
$data = @(
  @("Violación", 63, 295),
  @("Sodomia", 8, 0),
  @("Actos Lascivos", 151, 563),
  @("Incesto", 2, 3),
  @("Violación Técnica", 7, 27),
  @("Ley 54 (3.5)", 0, 2),
  @("Agresión Sexual", 6, 27),
  @("Maltrato", 133, 122),
  @("Pornografia infantil ", 13, 38),
  @("Hostigamiento Sexual", 0, 0),
  @("Maltrato Institucional", 7, 11),
  @("Trata Humana", 0, 0),
  @("Agresión", 1, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Column A width, matching the manually-resized width in the workbook
$ws.Columns.Item(1).ColumnWidth = 26.5

# Make the new sheet the active one with its own selection/cursor position
$ws.Select() | Out-Null
$ws.Range("E8").Select() | Out-Null
